$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to Text format
# so Excel stores them as text (matching the source inlineStr cells) instead of
# coercing them into numeric cells with floating-point artifacts.
$textCells = @("D4", "D5", "D6", "D7", "D12", "D16", "D18", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D37", "D38", "D43", "D44", "D45", "D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.113.63"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "3.908.77"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "487.20"
$ws.Range("E5").Value = "  +3.35%  "
$ws.Range("D6").Value = "146.52"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "42.87"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("D14").Value = "4.526.31"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "3.916.05"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "14.12"
$ws.Range("E16").Value = "  -6.70%  "
$ws.Range("D18").Value = "19.81"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "68.230.75"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "430.44"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "3.58"
$ws.Range("E22").Value = "  +5.35%  "
$ws.Range("E23").Value = "  +3.29%  "
$ws.Range("D24").Value = "87.33"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "11.37"
$ws.Range("E25").Value = "  +15.76%  "
$ws.Range("D26").Value = "11.35"
$ws.Range("E26").Value = "  +10.35%  "
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "38.11"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("D29").Value = "5.75"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "725.11"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").Value = "13.81"
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("D33").Value = "2.91"
$ws.Range("E33").Value = "  +5.71%  "
$ws.Range("D34").Value = "6.26"
$ws.Range("E34").Value = "  +16.76%  "
$ws.Range("D35").Value = "41.84"
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("D36").Value = "0.0₃0868"
$ws.Range("E36").Value = "  +4.22%  "
$ws.Range("D37").Value = "60.32"
$ws.Range("E37").Value = "  +3.95%  "
$ws.Range("D38").Value = "0.405"
$ws.Range("E38").Value = "  +19.60%  "
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  +16.17%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "3.15"
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("D44").Value = "2.91"
$ws.Range("E44").Value = "  +3.00%  "
$ws.Range("D45").Value = "0.140"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  +4.86%  "
$ws.Range("E48").Value = "  -4.23%  "
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "144.21"
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0338"
$ws.Range("E51").Value = "  +26.52%  "
